$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(57, 1).Value = "09.19.2022 09:39 (Kyiv+Israel) 06:39 (UTC) 15:39 (Japan) 12:09 (India)"
$ws.Cells.Item(57, 2).Value = 5.369
$ws.Cells.Item(57, 3).Value = -4.556
$ws.Cells.Item(57, 4).Value = "***"
$ws.Cells.Item(57, 5).Value = "***"
$ws.Cells.Item(58, 1).Value = "09.19.2022 10:39 (Kyiv+Israel) 07:39 (UTC) 16:39 (Japan) 13:09 (India)"
$ws.Cells.Item(58, 2).Value = 2.537
$ws.Cells.Item(58, 3).Value = -1.724
$ws.Cells.Item(58, 4).Value = "***"
$ws.Cells.Item(58, 5).Value = "***"
$ws.Cells.Item(59, 1).Value = "09.19.2022 10:43 (Kyiv+Israel) 07:43 (UTC) 16:43 (Japan) 13:13 (India)"
$ws.Cells.Item(59, 2).Value = "***"
$ws.Cells.Item(59, 3).Value = "***"
$ws.Cells.Item(59, 4).Value = 2.918
$ws.Cells.Item(59, 5).Value = -1.78
$ws.Cells.Item(60, 1).Value = "09.19.2022 11:34 (Kyiv+Israel) 08:34 (UTC) 17:34 (Japan) 14:04 (India)"
$ws.Cells.Item(60, 2).Value = "***"
$ws.Cells.Item(60, 3).Value = "***"
$ws.Cells.Item(60, 4).Value = 2.356
$ws.Cells.Item(60, 5).Value = -1.218
$ws.Cells.Item(61, 1).Value = "09.19.2022 11:53 (Kyiv+Israel) 08:53 (UTC) 17:53 (Japan) 14:23 (India)"
$ws.Cells.Item(61, 2).Value = 11.858
$ws.Cells.Item(61, 3).Value = -11.045
$ws.Cells.Item(61, 4).Value = "***"
$ws.Cells.Item(61, 5).Value = "***"
$ws.Cells.Item(62, 1).Value = "09.19.2022 13:44 (Kyiv+Israel) 10:44 (UTC) 19:44 (Japan) 16:14 (India)"
$ws.Cells.Item(62, 2).Value = 1.246
$ws.Cells.Item(62, 3).Value = -0.4330000000000001
$ws.Cells.Item(62, 4).Value = "***"
$ws.Cells.Item(62, 5).Value = "***"
$ws.Cells.Item(63, 1).Value = "09.19.2022 13:47 (Kyiv+Israel) 10:47 (UTC) 19:47 (Japan) 16:17 (India)"
$ws.Cells.Item(63, 2).Value = "***"
$ws.Cells.Item(63, 3).Value = "***"
$ws.Cells.Item(63, 4).Value = 1.142
$ws.Cells.Item(63, 5).Value = -0.004000000000000004
$ws.Cells.Item(64, 1).Value = "09.19.2022 15:05 (Kyiv+Israel) 12:05 (UTC) 21:05 (Japan) 17:35 (India)"
$ws.Cells.Item(64, 2).Value = 1.456
$ws.Cells.Item(64, 3).Value = -0.643
$ws.Cells.Item(64, 4).Value = "***"
$ws.Cells.Item(64, 5).Value = "***"
$ws.Cells.Item(65, 1).Value = "09.19.2022 15:22 (Kyiv+Israel) 12:22 (UTC) 21:22 (Japan) 17:52 (India)"
$ws.Cells.Item(65, 2).Value = 1.595
$ws.Cells.Item(65, 3).Value = -0.782
$ws.Cells.Item(65, 4).Value = "***"
$ws.Cells.Item(65, 5).Value = "***"

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(60, 1).Value = "09.19.2022 09:42 (Kyiv+Israel) 06:42 (UTC) 15:42 (Japan) 12:12 (India)"
$ws.Cells.Item(60, 2).Value = 1.114
$ws.Cells.Item(60, 3).Value = -0.4470000000000001
$ws.Cells.Item(60, 4).Value = "***"
$ws.Cells.Item(60, 5).Value = "***"
$ws.Cells.Item(61, 1).Value = "09.19.2022 09:45 (Kyiv+Israel) 06:45 (UTC) 15:45 (Japan) 12:15 (India)"
$ws.Cells.Item(61, 2).Value = "***"
$ws.Cells.Item(61, 3).Value = "***"
$ws.Cells.Item(61, 4).Value = 1.903
$ws.Cells.Item(61, 5).Value = -0.976
$ws.Cells.Item(62, 1).Value = "09.19.2022 11:37 (Kyiv+Israel) 08:37 (UTC) 17:37 (Japan) 14:07 (India)"
$ws.Cells.Item(62, 2).Value = "***"
$ws.Cells.Item(62, 3).Value = "***"
$ws.Cells.Item(62, 4).Value = 1.091
$ws.Cells.Item(62, 5).Value = -0.1639999999999999
$ws.Cells.Item(63, 1).Value = "09.19.2022 11:56 (Kyiv+Israel) 08:56 (UTC) 17:56 (Japan) 14:26 (India)"
$ws.Cells.Item(63, 2).Value = 1.301
$ws.Cells.Item(63, 3).Value = -0.6339999999999999
$ws.Cells.Item(63, 4).Value = "***"
$ws.Cells.Item(63, 5).Value = "***"
$ws.Cells.Item(64, 1).Value = "09.19.2022 15:07 (Kyiv+Israel) 12:07 (UTC) 21:07 (Japan) 17:37 (India)"
$ws.Cells.Item(64, 2).Value = 1.094
$ws.Cells.Item(64, 3).Value = -0.427
$ws.Cells.Item(64, 4).Value = "***"
$ws.Cells.Item(64, 5).Value = "***"
$ws.Cells.Item(65, 1).Value = "09.19.2022 15:25 (Kyiv+Israel) 12:25 (UTC) 21:25 (Japan) 17:55 (India)"
$ws.Cells.Item(65, 2).Value = 0.755
$ws.Cells.Item(65, 3).Value = -0.08799999999999997
$ws.Cells.Item(65, 4).Value = "***"
$ws.Cells.Item(65, 5).Value = "***"

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(53, 1).Value = "09.19.2022 09:44 (Kyiv+Israel) 06:44 (UTC) 15:44 (Japan) 12:14 (India)"
$ws.Cells.Item(53, 2).Value = 3.236
$ws.Cells.Item(53, 3).Value = -1.772
$ws.Cells.Item(53, 4).Value = "***"
$ws.Cells.Item(53, 5).Value = "***"
$ws.Cells.Item(54, 1).Value = "09.19.2022 09:47 (Kyiv+Israel) 06:47 (UTC) 15:47 (Japan) 12:17 (India)"
$ws.Cells.Item(54, 2).Value = "***"
$ws.Cells.Item(54, 3).Value = "***"
$ws.Cells.Item(54, 4).Value = 2.531
$ws.Cells.Item(54, 5).Value = -0.7270000000000001
$ws.Cells.Item(55, 1).Value = "09.19.2022 11:39 (Kyiv+Israel) 08:39 (UTC) 17:39 (Japan) 14:09 (India)"
$ws.Cells.Item(55, 2).Value = "***"
$ws.Cells.Item(55, 3).Value = "***"
$ws.Cells.Item(55, 4).Value = 2.021
$ws.Cells.Item(55, 5).Value = -0.2169999999999999
$ws.Cells.Item(56, 1).Value = "09.19.2022 11:59 (Kyiv+Israel) 08:59 (UTC) 17:59 (Japan) 14:29 (India)"
$ws.Cells.Item(56, 2).Value = 2.2
$ws.Cells.Item(56, 3).Value = -0.7360000000000002
$ws.Cells.Item(56, 4).Value = "***"
$ws.Cells.Item(56, 5).Value = "***"
$ws.Cells.Item(57, 1).Value = "09.19.2022 15:10 (Kyiv+Israel) 12:10 (UTC) 21:10 (Japan) 17:40 (India)"
$ws.Cells.Item(57, 2).Value = 3.657
$ws.Cells.Item(57, 3).Value = -2.193
$ws.Cells.Item(57, 4).Value = "***"
$ws.Cells.Item(57, 5).Value = "***"
$ws.Cells.Item(58, 1).Value = "09.19.2022 15:37 (Kyiv+Israel) 12:37 (UTC) 21:37 (Japan) 18:07 (India)"
$ws.Cells.Item(58, 2).Value = 1.64
$ws.Cells.Item(58, 3).Value = -0.1759999999999999
$ws.Cells.Item(58, 4).Value = "***"
$ws.Cells.Item(58, 5).Value = "***"

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(50, 1).Value = "09.19.2022 09:49 (Kyiv+Israel) 06:49 (UTC) 15:49 (Japan) 12:19 (India)"
$ws.Cells.Item(50, 2).Value = "***"
$ws.Cells.Item(50, 3).Value = "***"
$ws.Cells.Item(50, 4).Value = 1.356
$ws.Cells.Item(50, 5).Value = -0.4340000000000001
$ws.Cells.Item(51, 1).Value = "09.19.2022 09:52 (Kyiv+Israel) 06:52 (UTC) 15:52 (Japan) 12:22 (India)"
$ws.Cells.Item(51, 2).Value = 1.053
$ws.Cells.Item(51, 3).Value = -0.327
$ws.Cells.Item(51, 4).Value = "***"
$ws.Cells.Item(51, 5).Value = "***"
$ws.Cells.Item(52, 1).Value = "09.19.2022 11:42 (Kyiv+Israel) 08:42 (UTC) 17:42 (Japan) 14:12 (India)"
$ws.Cells.Item(52, 2).Value = "***"
$ws.Cells.Item(52, 3).Value = "***"
$ws.Cells.Item(52, 4).Value = 1.274
$ws.Cells.Item(52, 5).Value = -0.352
$ws.Cells.Item(53, 1).Value = "09.19.2022 12:09 (Kyiv+Israel) 09:09 (UTC) 18:09 (Japan) 14:39 (India)"
$ws.Cells.Item(53, 2).Value = 1.242
$ws.Cells.Item(53, 3).Value = -0.516
$ws.Cells.Item(53, 4).Value = "***"
$ws.Cells.Item(53, 5).Value = "***"
$ws.Cells.Item(54, 1).Value = "09.19.2022 15:18 (Kyiv+Israel) 12:18 (UTC) 21:18 (Japan) 17:48 (India)"
$ws.Cells.Item(54, 2).Value = 0.736
$ws.Cells.Item(54, 3).Value = -0.01000000000000001
$ws.Cells.Item(54, 4).Value = "***"
$ws.Cells.Item(54, 5).Value = "***"
$ws.Cells.Item(55, 1).Value = "09.19.2022 15:40 (Kyiv+Israel) 12:40 (UTC) 21:40 (Japan) 18:10 (India)"
$ws.Cells.Item(55, 2).Value = 0.836
$ws.Cells.Item(55, 3).Value = -0.11
$ws.Cells.Item(55, 4).Value = "***"
$ws.Cells.Item(55, 5).Value = "***"
